# "primer abono aparece en historial abonos"
#
# Bumps the invoice/order number from 119 to 121 (heading + concept
# line), advances the "Fecha" / "Fecha evento" header dates by one day
# (2019-02-26 -> 2019-02-27), and records the first abono payment that
# now shows up in the abonos history row:
#   fecha   2019-02-27 -> 2019-02-28
#   abono   39999      -> 30000
#   saldo   10000      -> 20000
#   total   49999      -> 50000
#
# NOTE on ordering: this runtime's Find.Execute always searches from the
# start of the document and acts on the first match, regardless of which
# Range/Cell it was invoked on - so instead of relying on range
# confinement, each replacement below targets text that is still unique
# (or whose first remaining occurrence is the correct one) at the moment
# it runs. In particular the abono-history date (old value "2019-02-27")
# is changed before the header dates are advanced to "2019-02-27",
# otherwise the new header text would collide with the search term.

$d = $word.ActiveDocument

# Order number: "No.  119" -> "No.  121" (title block)
$d.Content.Find.Execute("No.  119", $true, $false, $false, $false, $false, `
    $true, 0, $false, "No.  121", 1) | Out-Null

# Concept line: "Abono para compromiso con factura No: 119" -> "...121"
$d.Content.Find.Execute("factura No: 119", $true, $false, $false, $false, $false, `
    $true, 0, $false, "factura No: 121", 1) | Out-Null

# Abonos-history row date (still the only "2019-02-27" in the document
# at this point): 2019-02-27 -> 2019-02-28
$d.Content.Find.Execute("2019-02-27", $true, $false, $false, $false, $false, `
    $true, 0, $false, "2019-02-28", 1) | Out-Null

# Header "Fecha:" date: 2019-02-26 -> 2019-02-27
$d.Content.Find.Execute("2019-02-26", $true, $false, $false, $false, $false, `
    $true, 0, $false, "2019-02-27", 1) | Out-Null

# Header "Fecha evento:" date: 2019-02-26 -> 2019-02-27
$d.Content.Find.Execute("2019-02-26", $true, $false, $false, $false, $false, `
    $true, 0, $false, "2019-02-27", 1) | Out-Null

# Abonos-history amounts
$d.Content.Find.Execute("39999", $true, $false, $false, $false, $false, `
    $true, 0, $false, "30000", 1) | Out-Null

$d.Content.Find.Execute("10000", $true, $false, $false, $false, $false, `
    $true, 0, $false, "20000", 1) | Out-Null

$d.Content.Find.Execute("49999", $true, $false, $false, $false, $false, `
    $true, 0, $false, "50000", 1) | Out-Null
